$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A54").Value = 45986
$ws.Range("A54").Style = $ws.Range("A53").Style

$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = -0.08656168856399082
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = -0.1516437243033186
